# Add a new column U to the COVID19 deceased-cases time series sheet:
#   U1  = header text "06-10-2020" (same look as T1's header)
#   U2:U36 = the new day's numeric counts per state/UT

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell U1 --------------------------------------------------
# Give U1 a text value without Excel's automatic "looks like a date"
# re-interpretation (a leading apostrophe forces text entry), then copy
# only the formatting (font/border/alignment) from T1 - the last header
# cell - onto U1 so it visually matches the rest of the header row.
$ws.Range("U1").Value = "'06-10-2020"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

# --- Data cells U2:U36 -------------------------------------------------
$values = @{
    2  = 54
    3  = 6019
    4  = 19
    5  = 760
    6  = 924
    7  = 177
    8  = 1081
    9  = 2
    10 = 5542
    11 = 460
    12 = 3509
    13 = 1491
    14 = 224
    15 = 1252
    16 = 747
    17 = 9370
    18 = 859
    19 = 61
    20 = 2463
    21 = 38347
    22 = 75
    23 = 59
    24 = 0
    25 = 17
    26 = 924
    27 = 543
    28 = 3641
    29 = 1559
    30 = 46
    31 = 9846
    32 = 1181
    33 = 301
    34 = 669
    35 = 6092
    36 = 5255
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 21).Value = $values[$row]
}
